# Scheduled-runner refresh of the FFXIV Leve profitability workbook.
# Re-pulls market board "currentAveragePrice*" data (columns H-L) for the
# affected leve rows and recomputes the derived profit columns (M, N) on
# each of the eight job sheets. Columns A-G (leve metadata) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 530
$ws.Range("I5").Value = 487.0909
$ws.Range("J5").Value = 1002
$ws.Range("K5").Value = 487.0909
$ws.Range("L5").Value = 1002
$ws.Range("M5").Value = -372.0909
$ws.Range("N5").Value = -1232

$ws.Range("H19").Value = 4216.4375
$ws.Range("I19").Value = 608.6667
$ws.Range("J19").Value = 5049
$ws.Range("K19").Value = 608.6667
$ws.Range("L19").Value = 5049
$ws.Range("M19").Value = -433.6667
$ws.Range("N19").Value = -5399

$ws.Range("H76").Value = 9634.444
$ws.Range("I76").Value = 9344.200000000001
$ws.Range("J76").Value = 9997.25
$ws.Range("K76").Value = 9344.200000000001
$ws.Range("L76").Value = 9997.25
$ws.Range("M76").Value = -9029.200000000001
$ws.Range("N76").Value = -10627.25

$ws.Range("H79").Value = 9634.444
$ws.Range("I79").Value = 9344.200000000001
$ws.Range("J79").Value = 9997.25
$ws.Range("K79").Value = 9344.200000000001
$ws.Range("L79").Value = 9997.25
$ws.Range("M79").Value = -8252.200000000001
$ws.Range("N79").Value = -12181.25

$ws.Range("H137").Value = 1676.4878
$ws.Range("I137").Value = 1517.7646
$ws.Range("J137").Value = 2447.4285
$ws.Range("K137").Value = 4553.293799999999
$ws.Range("L137").Value = 7342.2855
$ws.Range("M137").Value = -2003.293799999999
$ws.Range("N137").Value = -12442.2855

$ws.Range("H139").Value = 93194.25
$ws.Range("J139").Value = 99259
$ws.Range("L139").Value = 99259
$ws.Range("N139").Value = -109539

$ws.Range("H140").Value = 78299.664
$ws.Range("J140").Value = 83959.60000000001
$ws.Range("L140").Value = 83959.60000000001
$ws.Range("N140").Value = -94319.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6077.916
$ws.Range("I32").Value = 4179.5234
$ws.Range("J32").Value = 24218.111
$ws.Range("K32").Value = 4179.5234
$ws.Range("L32").Value = 24218.111
$ws.Range("M32").Value = -3892.5234
$ws.Range("N32").Value = -24792.111

$ws.Range("H61").Value = 15929.931
$ws.Range("I61").Value = 2667.6482
$ws.Range("J61").Value = 55716.777
$ws.Range("K61").Value = 2667.6482
$ws.Range("L61").Value = 55716.777
$ws.Range("M61").Value = -2455.6482
$ws.Range("N61").Value = -56140.777

$ws.Range("H74").Value = 47458.31
$ws.Range("I74").Value = 26821.742
$ws.Range("K74").Value = 26821.742
$ws.Range("M74").Value = -25947.742

$ws.Range("H77").Value = 47458.31
$ws.Range("I77").Value = 26821.742
$ws.Range("K77").Value = 134108.71
$ws.Range("M77").Value = -129740.71

$ws.Range("H136").Value = 15929.931
$ws.Range("I136").Value = 2667.6482
$ws.Range("J136").Value = 55716.777
$ws.Range("K136").Value = 8002.944600000001
$ws.Range("L136").Value = 167150.331
$ws.Range("M136").Value = -5452.944600000001
$ws.Range("N136").Value = -172250.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2521.6086
$ws.Range("I20").Value = 2647.2778
$ws.Range("K20").Value = 2647.2778
$ws.Range("M20").Value = -2400.2778

$ws.Range("H99").Value = 3020.3333
$ws.Range("I99").Value = 1755.3636
$ws.Range("K99").Value = 1755.3636
$ws.Range("M99").Value = -257.3635999999999

$ws.Range("H134").Value = 3548.743
$ws.Range("I134").Value = 2420.44
$ws.Range("K134").Value = 7261.32
$ws.Range("M134").Value = -4726.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7656.391
$ws.Range("I25").Value = 4421
$ws.Range("J25").Value = 10145.154
$ws.Range("K25").Value = 4421
$ws.Range("L25").Value = 10145.154
$ws.Range("M25").Value = -4247
$ws.Range("N25").Value = -10493.154

$ws.Range("H31").Value = 2646.4482
$ws.Range("I31").Value = 2167.4167
$ws.Range("J31").Value = 2771.413
$ws.Range("K31").Value = 2167.4167
$ws.Range("L31").Value = 2771.413
$ws.Range("M31").Value = -1872.4167
$ws.Range("N31").Value = -3361.413

$ws.Range("H34").Value = 2646.4482
$ws.Range("I34").Value = 2167.4167
$ws.Range("J34").Value = 2771.413
$ws.Range("K34").Value = 2167.4167
$ws.Range("L34").Value = 2771.413
$ws.Range("M34").Value = -1965.4167
$ws.Range("N34").Value = -3175.413

$ws.Range("H41").Value = 22465.154
$ws.Range("J41").Value = 33331.332
$ws.Range("L41").Value = 33331.332
$ws.Range("N41").Value = -34187.332

$ws.Range("H99").Value = 3973.3684
$ws.Range("I99").Value = 3990.6296
$ws.Range("J99").Value = 3931
$ws.Range("K99").Value = 3990.6296
$ws.Range("L99").Value = 3931
$ws.Range("M99").Value = -2492.6296
$ws.Range("N99").Value = -6927

$ws.Range("H126").Value = 3973.3684
$ws.Range("I126").Value = 3990.6296
$ws.Range("J126").Value = 3931
$ws.Range("K126").Value = 11971.8888
$ws.Range("L126").Value = 11793
$ws.Range("M126").Value = -9501.888800000001
$ws.Range("N126").Value = -16733

$ws.Range("H141").Value = 91777.42999999999
$ws.Range("I141").Value = 24000
$ws.Range("J141").Value = 95166.3
$ws.Range("K141").Value = 24000
$ws.Range("L141").Value = 95166.3
$ws.Range("M141").Value = -18820
$ws.Range("N141").Value = -105526.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5203.6
$ws.Range("I80").Value = 3939.3333
$ws.Range("J80").Value = 7100
$ws.Range("K80").Value = 3939.3333
$ws.Range("L80").Value = 7100
$ws.Range("M80").Value = -2941.3333
$ws.Range("N80").Value = -9096

$ws.Range("H83").Value = 5203.6
$ws.Range("I83").Value = 3939.3333
$ws.Range("J83").Value = 7100
$ws.Range("K83").Value = 19696.6665
$ws.Range("L83").Value = 35500
$ws.Range("M83").Value = -14704.6665
$ws.Range("N83").Value = -45484

$ws.Range("H126").Value = 28824.5
$ws.Range("I126").Value = 32262.5
$ws.Range("K126").Value = 96787.5
$ws.Range("M126").Value = -94317.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1019.85
$ws.Range("I16").Value = 790.9167
$ws.Range("J16").Value = 1363.25
$ws.Range("K16").Value = 790.9167
$ws.Range("L16").Value = 1363.25
$ws.Range("M16").Value = -620.9167
$ws.Range("N16").Value = -1703.25

$ws.Range("H40").Value = 3472.239
$ws.Range("I40").Value = 3244.6843
$ws.Range("J40").Value = 4553.125
$ws.Range("K40").Value = 3244.6843
$ws.Range("L40").Value = 4553.125
$ws.Range("M40").Value = -3108.6843
$ws.Range("N40").Value = -4825.125

$ws.Range("H46").Value = 2308.2964
$ws.Range("I46").Value = 2217
$ws.Range("J46").Value = 2326.5557
$ws.Range("K46").Value = 2217
$ws.Range("L46").Value = 2326.5557
$ws.Range("M46").Value = -2029
$ws.Range("N46").Value = -2702.5557

$ws.Range("H68").Value = 2688.5557
$ws.Range("I68").Value = 2633.6667
$ws.Range("J68").Value = 2716
$ws.Range("K68").Value = 2633.6667
$ws.Range("L68").Value = 2716
$ws.Range("M68").Value = -1884.6667
$ws.Range("N68").Value = -4214

$ws.Range("H71").Value = 2688.5557
$ws.Range("I71").Value = 2633.6667
$ws.Range("J71").Value = 2716
$ws.Range("K71").Value = 13168.3335
$ws.Range("L71").Value = 13580
$ws.Range("M71").Value = -9424.333500000001
$ws.Range("N71").Value = -21068

$ws.Range("H82").Value = 2273.7
$ws.Range("J82").Value = 3148.75
$ws.Range("L82").Value = 3148.75
$ws.Range("N82").Value = -3870.75

$ws.Range("H85").Value = 2273.7
$ws.Range("J85").Value = 3148.75
$ws.Range("L85").Value = 3148.75
$ws.Range("N85").Value = -5644.75

$ws.Range("H122").Value = 3156.5208
$ws.Range("I122").Value = 2941.4773
$ws.Range("J122").Value = 5522
$ws.Range("K122").Value = 8824.4319
$ws.Range("L122").Value = 16566
$ws.Range("M122").Value = -6374.4319
$ws.Range("N122").Value = -21466

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 17647500
$ws.Range("I5").Value = 17647500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 17647500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -17647388

$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1939

$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -9696

$ws.Range("H126").Value = 2013.3529
$ws.Range("I126").Value = 2013.3529
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6040.0587
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3570.0587

$ws.Range("H132").Value = 4002.4
$ws.Range("I132").Value = 3028.7368
$ws.Range("J132").Value = 22502
$ws.Range("K132").Value = 9086.2104
$ws.Range("L132").Value = 67506
$ws.Range("M132").Value = -6556.2104
$ws.Range("N132").Value = -72566

$ws.Range("H136").Value = 3388.8245
$ws.Range("I136").Value = 3140.608
$ws.Range("J136").Value = 5498.6665
$ws.Range("K136").Value = 9421.824000000001
$ws.Range("L136").Value = 16495.9995
$ws.Range("M136").Value = -6871.824000000001
$ws.Range("N136").Value = -21595.9995

$ws.Range("H141").Value = 164500
$ws.Range("J141").Value = 164500
$ws.Range("L141").Value = 164500
$ws.Range("N141").Value = -174860

# Special cases: WVR rows where the N column (LeveProfitHQ) cell
# is removed entirely (formula result became blank/N-A upstream)
$wsWvr = $wb.Worksheets.Item("WVR")
$wsWvr.Range("N5").ClearContents()
$wsWvr.Range("N126").ClearContents()

# Special cases: WVR rows where the N column cell is newly added
$wsWvr.Range("N81").Value = -5122
$wsWvr.Range("N84").Value = -25608
